$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F5").Value = -1
$ws.Range("F8").Value = 4
$ws.Range("F9").Value = -8
$ws.Range("F12").Value = -1
$ws.Range("F15").Value = -5
$ws.Range("F16").Value = -1
$ws.Range("F18").Value = 3
$ws.Range("F19").Value = 1
$ws.Range("F20").Value = 0
$ws.Range("F21").Value = -1
$ws.Range("F24").Value = 2
$ws.Range("F25").Value = 2
$ws.Range("F26").Value = 0
$ws.Range("F27").Value = -3
$ws.Range("F31").Value = -1
$ws.Range("F34").Value = -1
$ws.Range("F42").Value = 1
$ws.Range("F47").Value = -4
$ws.Range("F48").Value = 4
$ws.Range("F53").Value = -7
$ws.Range("F54").Value = -2
$ws.Range("F57").Value = -7
$ws.Range("F59").Value = 3
